$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.272.67"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.966.74"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'356.38"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "'109.28"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("D7").Value = "'0.569"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.628"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").Value = "'38.93"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "'0.0874"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "'19.35"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").Value = "3.439.01"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "'7.77"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "3.002.10"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "'0.984"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "52.265.29"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'3.51"
$ws.Range("E19").Value = "  +6.09%  "
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "'13.81"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Value = "'70.36"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'268.49"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'2.80"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "'0.177"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("D27").Value = "'27.16"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").Value = "'7.64"
$ws.Range("E28").Value = "  +15.81%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'0.108"
$ws.Range("E30").Value = "  +4.21%  "
$ws.Range("D31").Value = "'10.41"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "'37.33"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").Value = "'52.45"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").Value = "'0.0442"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'3.22"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "'18.04"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").Value = "'119.20"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("E46").Value = "  -4.59%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.139.64"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'3.44"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("E50").Value = "  -7.92%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.37"
$ws.Range("E51").Value = "  +1.72%  "
